# LoginData4.xlsx - 3rd commit
# Update the "Automation*" test-data row (row 2) on Sheet1 and remove the
# now-unused second data row (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: rename Rakesh-prefixed model/asset/supplier values to the
# Automation* equivalents, and swap the default-location value.
$ws.Range("D2").Value = "AutomationModel"
$ws.Range("F2").Value = "AutomationAsset"
$ws.Range("H2").Value = "AutomationSupplier"
$ws.Range("M2").Value = "parola"

# Row 3 (Rakesh22 / us-9877 / ASAN / ...) is no longer needed.
$ws.Rows.Item(3).Delete()

# The active selection moved to C5.
$ws.Range("C5").Select()

# Widen column D (Model) and give column H (Supplier) an explicit width
# matching column E's.
$ws.Columns.Item(4).ColumnWidth = 19
$ws.Columns.Item(8).ColumnWidth = 20.833333333333332
